$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Styles: give C1/D1 (and their computational-comparison counterparts) new
# border-only styles (top+bottom, and top+right+bottom) cloned from the base
# "Normal" style, matching the new cellXfs entries 2 and 3 in the diff.

# Sheet "quality_comparison": C1 -> top+bottom border, D1 -> top+right+bottom border
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1           # all four edges thin ...
$c1.Borders.Item(7).LineStyle = -4142   # ... then drop left
$c1.Borders.Item(10).LineStyle = -4142  # ... and drop right

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1           # all four edges thin ...
$d1.Borders.Item(7).LineStyle = -4142   # ... then drop left only

# Propagate the exact same styles (by format copy, so no extra style entries
# are created) to the matching cells in "computational_comparison".
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Content edits ---

# Rename "fedcore" header to "approach" (anonymized) on both sheets.
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Fix the -0 -> 0 change value on the quality sheet.
$ws1.Range("D4").Value = 0

# Remove the stray empty inline-string cell G5 on the computational sheet.
$ws2.Range("G5").ClearContents()

Write-Output "done"
